# Add a new worksheet "ODI Batting Extra" at the end of the workbook,
# containing per-match batting extras for the player.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it is appended at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# --- Header row (bold / centered / bordered, matching the other sheets) ---
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Copy the header formatting used by the existing sheets (bold, centered,
# bordered) so the new header row re-uses the same style instead of creating
# a near-duplicate one.
$headerFormat = $wb.Worksheets.Item("Player Info").Range("A1")
$headerFormat.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# --- Data row ---
# MATCH_CODE ("4484") and the NUM_4 / NUM_6 counters ("0") are stored as
# text, matching the look-alike-numeric text values used elsewhere in this
# workbook, so force them to text with a leading apostrophe.
$ws.Range("A2").Value = "'4484"
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = "'0"
$ws.Range("D2").Value = "'0"
$ws.Range("E2").Value = "'"
$ws.Range("F2").Value = "NO"

# Keep the originally active sheet ("Player Info") selected, since adding a
# new worksheet makes it the active one by default.
$wb.Worksheets.Item(1).Activate()
